$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = -13.376
$ws.Range("B9").Value = 5.567000000000001
$ws.Range("C9").Value = -11.18
$ws.Range("D9").Value = -7.173
$ws.Range("C11").Value = -12.464
$ws.Range("B18").Value = 5.275
$ws.Range("B20").Value = 6.601999999999999
$ws.Range("C23").Value = -13.004
$ws.Range("C24").Value = -12.45
$ws.Range("C26").Value = -11.576
$ws.Range("B27").Value = 5.81
$ws.Range("D27").Value = -7.731999999999999
$ws.Range("D29").Value = -7.538999999999999
$ws.Range("D32").Value = -7.353999999999999
$ws.Range("C34").Value = -12.098
$ws.Range("B35").Value = 8.340999999999999
$ws.Range("C35").Value = -12.547
$ws.Range("D37").Value = -7.87
$ws.Range("D38").Value = -7.449
$ws.Range("D41").Value = -8.106999999999999
$ws.Range("D45").Value = -7.525
$ws.Range("C48").Value = -11.689
$ws.Range("C49").Value = -13.088
$ws.Range("D51").Value = -8.010000000000002
$ws.Range("C52").Value = -11.662
$ws.Range("D57").Value = -8.296000000000001
$ws.Range("D64").Value = -7.888
$ws.Range("C66").Value = -11.187
$ws.Range("C67").Value = -11.492
$ws.Range("B69").Value = 5.930999999999999
$ws.Range("B76").Value = 6.11
$ws.Range("B78").Value = 8.75
$ws.Range("C78").Value = -12.059
$ws.Range("C80").Value = -12.048
$ws.Range("B82").Value = 5.211
$ws.Range("D82").Value = -8.181000000000001
$ws.Range("B83").Value = 5.598000000000001
$ws.Range("B93").Value = 5.008
$ws.Range("D93").Value = -7.064
$ws.Range("C99").Value = -11.823
$ws.Range("D102").Value = -7.709999999999999
$ws.Range("C104").Value = -13.099
$ws.Range("D105").Value = -7.876
